$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.8
$ws.Range("G2").Value = 3.85
$ws.Range("H2").Value = 2.54
$ws.Range("I2").Value = 3.15
$ws.Range("K2").Value = 3.35
$ws.Range("L2").Value = 1.56
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 2.14
$ws.Range("O2").Value = 1.66
$ws.Range("P2").Value = 1.34
$ws.Range("Q2").Value = 2.72
$ws.Range("S2").Value = 6.8
$ws.Range("T2").Value = 2.28
$ws.Range("U2").Value = 1.64
$ws.Range("V2").Value = 1.46
$ws.Range("W2").Value = 1.4
$ws.Range("F3").Value = 1.17
$ws.Range("H3").Value = 16
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 11
$ws.Range("S3").Value = 1.73
$ws.Range("U3").Value = 2
$ws.Range("W3").Value = 5.8
$ws.Range("X3").Value = 70
$ws.Range("AB3").Value = 19
$ws.Range("AC3").Value = 28
$ws.Range("AG3").Value = 16
$ws.Range("AH3").Value = 44
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 16.5
$ws.Range("AL3").Value = 42
$ws.Range("AM3").Value = 180
$ws.Range("AN3").Value = 2.76
$ws.Range("F4").Value = 1.63
$ws.Range("G4").Value = 1000
$ws.Range("I4").Value = 8.6
$ws.Range("K4").Value = 5.2
$ws.Range("L4").Value = 1.45
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 1.45
$ws.Range("Q4").Value = 2.24
$ws.Range("R4").Value = 1.13
$ws.Range("S4").Value = 2.66
$ws.Range("V4").Value = 1.13
$ws.Range("W4").Value = 1.01
$ws.Range("R5").Value = 1.24
$ws.Range("X5").Value = 9.800000000000001
$ws.Range("F6").Value = 1.63
$ws.Range("G6").Value = 1.71
$ws.Range("H6").Value = 4.3
$ws.Range("I6").Value = 5.1
$ws.Range("R6").Value = 1.9
$ws.Range("S6").Value = 1.94
$ws.Range("V6").Value = 1.25
$ws.Range("W6").Value = 2.32
$ws.Range("Y6").Value = 34
$ws.Range("AB6").Value = 18.5
$ws.Range("AI6").Value = 42
$ws.Range("AJ6").Value = 21
$ws.Range("AN6").Value = 5.5
$ws.Range("AO6").Value = 29
$ws.Range("X7").Value = 17.5
$ws.Range("AN7").Value = 110
$ws.Range("AO7").Value = 13.5
$ws.Range("G8").Value = 1.77
$ws.Range("I8").Value = 8.199999999999999
$ws.Range("J8").Value = 3.15
$ws.Range("L8").Value = 1.43
$ws.Range("Q8").Value = 2.2
$ws.Range("S8").Value = 3.5
$ws.Range("V8").Value = 1.13
$ws.Range("W8").Value = 2.28
$ws.Range("F9").Value = 1.45
$ws.Range("O9").Value = 1.38
$ws.Range("AA9").Value = 480
$ws.Range("AL9").Value = 50
$ws.Range("AN9").Value = 8.800000000000001
$ws.Range("F10").Value = 1.64
$ws.Range("G10").Value = 1.65
$ws.Range("K10").Value = 4.5
$ws.Range("M10").Value = 1.06
$ws.Range("P10").Value = 2.14
$ws.Range("Q10").Value = 1.84
$ws.Range("T10").Value = 1.89
$ws.Range("W10").Value = 2.54
$ws.Range("AA10").Value = 160
$ws.Range("AG10").Value = 9.6
$ws.Range("F11").Value = 2.04
$ws.Range("G11").Value = 2.56
$ws.Range("H11").Value = 3.45
$ws.Range("K11").Value = 4.8
$ws.Range("N11").Value = 1.01
$ws.Range("P11").Value = 1.25
$ws.Range("Q11").Value = 1.02
$ws.Range("T11").Value = 1.03
$ws.Range("U11").Value = 2
$ws.Range("W11").Value = 1.64
